$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "MCH152-1"
$ws.Range("C2").Value = "NETHEERLANDS, OMROEP VOOR, RADIO FREEDOM, LEAFLETS, MAGAZINES, ARTICLES, STICKERS, POSTCARDS, BOOKS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 21N | GRAP COUNT NUMER: NONE"

# Row 3
$ws.Range("A3").Value = "MCH152-2"
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: 21N | GRAP COUNT NUMER: NONE"
